# Auto-generated edit script: updates modification-driver trajectories
# for sheet "strategy_id-0" (mexico, calibrated).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Row 97: J97:AR97 (35 cells)
$vals = New-Object "object[,]" 1,35
$vals[0,0] = 1
$vals[0,1] = 1
$vals[0,2] = 1
$vals[0,3] = 1
$vals[0,4] = 1
$vals[0,5] = 1
$vals[0,6] = 1
$vals[0,7] = 1
$vals[0,8] = 1
$vals[0,9] = 1
$vals[0,10] = 1
$vals[0,11] = 1
$vals[0,12] = 1
$vals[0,13] = 1
$vals[0,14] = 1
$vals[0,15] = 1
$vals[0,16] = 1
$vals[0,17] = 1
$vals[0,18] = 1
$vals[0,19] = 1
$vals[0,20] = 1
$vals[0,21] = 1
$vals[0,22] = 1
$vals[0,23] = 1
$vals[0,24] = 1
$vals[0,25] = 1
$vals[0,26] = 1
$vals[0,27] = 1
$vals[0,28] = 1
$vals[0,29] = 1
$vals[0,30] = 1
$vals[0,31] = 1
$vals[0,32] = 1
$vals[0,33] = 1
$vals[0,34] = 1
$ws.Range("J97:AR97").Value = $vals

# Row 110: J110:AR110 (35 cells)
$vals = New-Object "object[,]" 1,35
$vals[0,0] = 1
$vals[0,1] = 1
$vals[0,2] = 1
$vals[0,3] = 1
$vals[0,4] = 1
$vals[0,5] = 1
$vals[0,6] = 1
$vals[0,7] = 1
$vals[0,8] = 1
$vals[0,9] = 1
$vals[0,10] = 1
$vals[0,11] = 1
$vals[0,12] = 1
$vals[0,13] = 1
$vals[0,14] = 1
$vals[0,15] = 1
$vals[0,16] = 1
$vals[0,17] = 1
$vals[0,18] = 1
$vals[0,19] = 1
$vals[0,20] = 1
$vals[0,21] = 1
$vals[0,22] = 1
$vals[0,23] = 1
$vals[0,24] = 1
$vals[0,25] = 1
$vals[0,26] = 1
$vals[0,27] = 1
$vals[0,28] = 1
$vals[0,29] = 1
$vals[0,30] = 1
$vals[0,31] = 1
$vals[0,32] = 1
$vals[0,33] = 1
$vals[0,34] = 1
$ws.Range("J110:AR110").Value = $vals

# Row 132: J132:AS132 (36 cells)
$vals = New-Object "object[,]" 1,36
$vals[0,0] = 34994250
$vals[0,1] = 36028600
$vals[0,2] = 36511800
$vals[0,3] = 36489150
$vals[0,4] = 34095800
$vals[0,5] = 36111650
$vals[0,6] = 36833883
$vals[0,7] = 37570560.66
$vals[0,8] = 38321971.8732
$vals[0,9] = 39088411.30855
$vals[0,10] = 39870179.5391
$vals[0,11] = 40667583.12505
$vals[0,12] = 41480934.7874
$vals[0,13] = 42310553.48345
$vals[0,14] = 43156764.5578
$vals[0,15] = 44019899.84805
$vals[0,16] = 44900297.84335
$vals[0,17] = 45798303.79765
$vals[0,18] = 46714269.87315001
$vals[0,19] = 47648555.2762
$vals[0,20] = 48601526.3781
$vals[0,21] = 49573556.90385
$vals[0,22] = 50565028.0454
$vals[0,23] = 51576328.6051
$vals[0,24] = 52607855.1769
$vals[0,25] = 53660012.28225
$vals[0,26] = 54733212.52865
$vals[0,27] = 55827876.77575
$vals[0,28] = 56944434.31655
$vals[0,29] = 58083322.9982
$vals[0,30] = 59244989.4636
$vals[0,31] = 60429889.24955
$vals[0,32] = 61638487.0359
$vals[0,33] = 62871256.7739
$vals[0,34] = 64128681.9127
$vals[0,35] = 65411255.55035
$ws.Range("J132:AS132").Value = $vals

# Row 134: J134:AS134 (36 cells)
$vals = New-Object "object[,]" 1,36
$vals[0,0] = 268117.585026178
$vals[0,1] = 248575.329928616
$vals[0,2] = 274810.17051554
$vals[0,3] = 278291.884002256
$vals[0,4] = 292482.3242657
$vals[0,5] = 259930.906759628
$vals[0,6] = 259931
$vals[0,7] = 261316.082362874
$vals[0,8] = 262708.545350422
$vals[0,9] = 264108.428291436
$vals[0,10] = 265515.770724286
$vals[0,11] = 266930.612398018
$vals[0,12] = 268352.993273496
$vals[0,13] = 269782.953524514
$vals[0,14] = 271220.533538944
$vals[0,15] = 272665.77391987
$vals[0,16] = 274118.71548673
$vals[0,17] = 275579.399276484
$vals[0,18] = 277047.866544756
$vals[0,19] = 278524.158767012
$vals[0,20] = 280008.317639722
$vals[0,21] = 281500.385081546
$vals[0,22] = 283000.403234512
$vals[0,23] = 284508.41446521
$vals[0,24] = 286024.461365982
$vals[0,25] = 287548.586756138
$vals[0,26] = 289080.83368315
$vals[0,27] = 290621.24542388
$vals[0,28] = 292169.865485792
$vals[0,29] = 293726.737608194
$vals[0,30] = 295291.905763456
$vals[0,31] = 296865.41415827
$vals[0,32] = 298447.307234886
$vals[0,33] = 300037.629672372
$vals[0,34] = 301636.426387876
$vals[0,35] = 303243.742537894
$ws.Range("J134:AS134").Value = $vals

# Row 135: J135:AS135 (36 cells)
$vals = New-Object "object[,]" 1,36
$vals[0,0] = 1091477.476575034
$vals[0,1] = 1025014.0815695
$vals[0,2] = 1009615.37098538
$vals[0,3] = 1035600.95417141
$vals[0,4] = 1056164.473529098
$vals[0,5] = 929161.87886073
$vals[0,6] = 929161.8
$vals[0,7] = 908485.999546466
$vals[0,8] = 888270.279053596
$vals[0,9] = 868504.40077651
$vals[0,10] = 849178.354781644
$vals[0,11] = 830282.353877468
$vals[0,12] = 811806.828658006
$vals[0,13] = 793742.4226566521
$vals[0,14] = 776079.9876078221
$vals[0,15] = 758810.578814046
$vals[0,16] = 741925.450616146
$vals[0,17] = 725416.051964224
$vals[0,18] = 709274.022087186
$vals[0,19] = 693491.186258646
$vals[0,20] = 678059.551657042
$vals[0,21] = 662971.303317868
$vals[0,22] = 648218.8001759819
$vals[0,23] = 633794.571195984
$vals[0,24] = 619691.3115886881
$vals[0,25] = 605901.879111808
$vals[0,26] = 592419.290452936
$vals[0,27] = 579236.7176930259
$vals[0,28] = 566347.484848562
$vals[0,29] = 553745.064490676
$vals[0,30] = 541423.074439492
$vals[0,31] = 529375.274532038
$vals[0,32] = 517595.563462062
$vals[0,33] = 506077.975690184
$vals[0,34] = 494816.678422798
$vals[0,35] = 483805.968658202
$ws.Range("J135:AS135").Value = $vals

# Row 136: J136:AS136 (36 cells)
$vals = New-Object "object[,]" 1,36
$vals[0,0] = 37016.0726055204
$vals[0,1] = 51714.320922467
$vals[0,2] = 41151.1357983976
$vals[0,3] = 44005.6200703732
$vals[0,4] = 7785.49232923952
$vals[0,5] = 42065.8943532988
$vals[0,6] = 42065.9
$vals[0,7] = 131670.5195680288
$vals[0,8] = 412142.037215764
$vals[0,9] = 1290046.241160312
$vals[0,10] = 4037975.14947644
$vals[0,11] = 12639270.42888302
$vals[0,12] = 39562194.1841694
$vals[0,13] = 123833667.2573466
$vals[0,14] = 387611897.232418
$vals[0,15] = 1213264423.186998
$vals[0,16] = 3797640297.11564
$vals[0,17] = 11886998044.8225
$vals[0,18] = 37207505572.5878
$vals[0,19] = 116463253860.5606
$vals[0,20] = 364541758203.278
$vals[0,21] = 1141052555796.22
$vals[0,22] = 3571609852068.18
$vals[0,23] = 11179499901728.2
$vals[0,24] = 34992964861592.8
$vals[0,25] = 109531517560582.4
$vals[0,26] = 342844722834328
$vals[0,27] = 1073138641673010
$vals[0,28] = 3359032435241226
$vals[0,29] = 10514111096971020.0
$vals[0,30] = 32910230636552400.0
$vals[0,31] = 103012348886354800.0
$vals[0,32] = 322439065841676288.0
$vals[0,33] = 1009266872416930048.0
$vals[0,34] = 3159107340480924160.0
$vals[0,35] = 9888325339343665152.0
$ws.Range("J136:AS136").Value = $vals

# Row 138: J138:AS138 (36 cells)
$vals = New-Object "object[,]" 1,36
$vals[0,0] = 226859846.2232036
$vals[0,1] = 190618182.9631508
$vals[0,2] = 158908712.3046357
$vals[0,3] = 168557916.4624421
$vals[0,4] = 164321605.9565982
$vals[0,5] = 153095485.5635689
$vals[0,6] = 153095470
$vals[0,7] = 164853622.5114702
$vals[0,8] = 177514833.4248839
$vals[0,9] = 191148460.1053995
$vals[0,10] = 205829186.7542801
$vals[0,11] = 221637433.5245389
$vals[0,12] = 238659797.0577788
$vals[0,13] = 256989524.8554854
$vals[0,14] = 276727026.083321
$vals[0,15] = 297980421.6066062
$vals[0,16] = 320866136.2700299
$vals[0,17] = 345509536.6660641
$vals[0,18] = 372045617.8857573
$vals[0,19] = 400619743.0138541
$vals[0,20] = 431388439.4191945
$vals[0,21] = 464520256.2023826
$vals[0,22] = 500196687.4978007
$vals[0,23] = 538613166.6877556
$vals[0,24] = 579980136.9749916
$vals[0,25] = 624524204.1781259
$vals[0,26] = 672489378.0649232
$vals[0,27] = 724138409.0233898
$vals[0,28] = 779754227.3928716
$vals[0,29] = 839641493.3397282
$vals[0,30] = 904128265.7677237
$vals[0,31] = 973567799.4053155
$vals[0,32] = 1048340479.914178
$vals[0,33] = 1128855907.619374
$vals[0,34] = 1215555141.27575
$vals[0,35] = 1308913114.161708
$ws.Range("J138:AS138").Value = $vals

